# Commit: "Fix: Remove a slide."
#
# The deck had a duplicated "刪除記錄功能" (delete-record feature) slide at
# positions 5 and 6. The author removed the duplicate that sat at slide
# position 5 (sldId 258), so the remaining slides shift up: the slide that
# used to be #6 (sldId 262) becomes #5, and the old #7 (sldId 259) becomes #6.
$p = $ppt.ActivePresentation
$p.Slides.Item(5).Delete()
